$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename columns to snake_case machine-readable names ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Normalize capitalization of connector words (de/del/la/las/los/y/el) ---
# in municipio / estado names to title case, per upstream cleaning script fix ---
$ws.Range('B7').Value = 'Pabellón De Arteaga'
$ws.Range('B8').Value = 'Rincón De Romos'
$ws.Range('B12').Value = 'Playas De Rosarito'
$ws.Range('B46').Value = 'San Cristóbal De Las Casas'
$ws.Range('B74').Value = 'Coyame Del Sotol'
$ws.Range('B85').Value = 'Guadalupe Y Calvo'
$ws.Range('B88').Value = 'Hidalgo Del Parral'
$ws.Range('B113').Value = 'San Francisco De Borja'
$ws.Range('B114').Value = 'San Francisco De Conchos'
$ws.Range('B115').Value = 'San Francisco Del Oro'
$ws.Range('B123').Value = 'Valle De Zaragoza'
$ws.Range('B139').Value = 'San Juan De Sabinas'
$ws.Range('A150').Value = 'Ciudad De México'
$ws.Range('B154').Value = 'Cuajimalpa De Morelos'
$ws.Range('B169').Value = 'Coneto De Comonfort'
$ws.Range('B183').Value = 'Nombre De Dios'
$ws.Range('B187').Value = 'Pánuco De Coronado'
$ws.Range('B194').Value = 'San Juan De Guadalupe'
$ws.Range('B195').Value = 'San Juan Del Río'
$ws.Range('B196').Value = 'San Luis Del Cordero'
$ws.Range('B197').Value = 'San Pedro Del Gallo'
$ws.Range('A207').Value = 'Estado De México'
$ws.Range('B209').Value = 'Almoloya De Juárez'
$ws.Range('B213').Value = 'Atizapán De Zaragoza'
$ws.Range('B217').Value = 'Coacalco De Berriozábal'
$ws.Range('B220').Value = 'Ecatepec De Morelos'
$ws.Range('B223').Value = 'Ixtapan De La Sal'
$ws.Range('B228').Value = 'Naucalpan De Juárez'
$ws.Range('B231').Value = 'San Felipe Del Progreso'
$ws.Range('B232').Value = 'San Martín De Las Pirámides'
$ws.Range('B238').Value = 'Tenango Del Valle'
$ws.Range('B242').Value = 'Tlalnepantla De Baz'
$ws.Range('B247').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B248').Value = 'Villa Del Carbón'
$ws.Range('B255').Value = 'San Miguel De Allende'
$ws.Range('B256').Value = 'Apaseo El Alto'
$ws.Range('B263').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B267').Value = 'Jaral Del Progreso'
$ws.Range('B274').Value = 'Purísima Del Rincón'
$ws.Range('B279').Value = 'San Francisco Del Rincón'
$ws.Range('B281').Value = 'San Luis De La Paz'
$ws.Range('B283').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B285').Value = 'Silao De La Victoria'
$ws.Range('B290').Value = 'Valle De Santiago'
$ws.Range('B295').Value = 'Acapulco De Juárez'
$ws.Range('B296').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B299').Value = 'Atenango Del Río'
$ws.Range('B300').Value = 'Atoyac De Álvarez'
$ws.Range('B301').Value = 'Ayutla De Los Libres'
$ws.Range('B303').Value = 'Buenavista De Cuéllar'
$ws.Range('B304').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B305').Value = 'Coahuayutla De José María Izazaga'
$ws.Range('B308').Value = 'Coyuca De Benítez'
$ws.Range('B309').Value = 'Coyuca De Catalán'
$ws.Range('B312').Value = 'Cuetzala Del Progreso'
$ws.Range('B313').Value = 'Cutzamala De Pinzón'
$ws.Range('B317').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B318').Value = 'Iguala De La Independencia'
$ws.Range('B319').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B329').Value = 'Taxco De Alarcón'
$ws.Range('B331').Value = 'Técpan De Galeana'
$ws.Range('B333').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B334').Value = 'Tlalixtaquilla De Maldonado'
$ws.Range('B339').Value = 'Agua Blanca De Iturbide'
$ws.Range('B346').Value = 'Cuautepec De Hinojosa'
$ws.Range('B349').Value = 'Huasca De Ocampo'
$ws.Range('B355').Value = 'Mineral Del Monte'
$ws.Range('B356').Value = 'Mixquiahuala De Juárez'
$ws.Range('B357').Value = 'Nopala De Villagrán'
$ws.Range('B358').Value = 'Pachuca De Soto'
$ws.Range('B359').Value = 'Progreso De Obregón'
$ws.Range('B366').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B367').Value = 'Tula De Allende'
$ws.Range('B368').Value = 'Tulancingo De Bravo'
$ws.Range('B373').Value = 'Acatlán De Juárez'
$ws.Range('B378').Value = 'Atotonilco El Alto'
$ws.Range('B380').Value = 'Autlán De Navarro'
$ws.Range('B389').Value = 'Concepción De Buenos Aires'
$ws.Range('B395').Value = 'Encarnación De Díaz'
$ws.Range('B399').Value = 'Huejuquilla El Alto'
$ws.Range('B400').Value = 'Ixtlahuacán Del Río'
$ws.Range('B404').Value = 'Jilotlán De Los Dolores'
$ws.Range('B408').Value = 'Lagos De Moreno'
$ws.Range('B413').Value = 'Ojuelos De Jalisco'
$ws.Range('B417').Value = 'San Diego De Alejandría'
$ws.Range('B419').Value = 'San Juan De Los Lagos'
$ws.Range('B420').Value = 'San Juanito De Escobedo'
$ws.Range('B422').Value = 'San Miguel El Alto'
$ws.Range('B424').Value = 'Tamazula De Gordiano'
$ws.Range('B427').Value = 'Teocuitatlán De Corona'
$ws.Range('B428').Value = 'Tepatitlán De Morelos'
$ws.Range('B430').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B436').Value = 'Unión De Tula'
$ws.Range('B440').Value = 'Yahualica De González Gallo'
$ws.Range('B442').Value = 'Zapotitlán De Vadillo'
$ws.Range('B443').Value = 'Zapotlán El Grande'
$ws.Range('B453').Value = 'Coalcomán De Vázquez Pallares'
$ws.Range('B507').Value = 'Coatlán Del Río'
$ws.Range('B513').Value = 'Puente De Ixtla'
$ws.Range('B517').Value = 'Tetela Del Volcán'
$ws.Range('B531').Value = 'Santa María Del Oro'
$ws.Range('B548').Value = 'Mier Y Noriega'
$ws.Range('B550').Value = 'San Nicolás De Los Garza'
$ws.Range('B556').Value = 'Guevea De Humboldt'
$ws.Range('B557').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B558').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B559').Value = 'Ixtlán De Juárez'
$ws.Range('B560').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B563').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B564').Value = 'Oaxaca De Juárez'
$ws.Range('B565').Value = 'Ocotlán De Morelos'
$ws.Range('B566').Value = 'Putla Villa De Guerrero'
$ws.Range('B567').Value = 'San Antonio De La Cal'
$ws.Range('B596').Value = 'Tamazulápam Del Espíritu Santo'
$ws.Range('B597').Value = 'Teotitlán De Flores Magón'
$ws.Range('B598').Value = 'Teotitlán Del Valle'
$ws.Range('B599').Value = 'Tepelmeme Villa De Morelos'
$ws.Range('B601').Value = 'Tlacolula De Matamoros'
$ws.Range('B602').Value = 'Villa De Etla'
$ws.Range('B603').Value = 'Villa Sola De Vega'
$ws.Range('B604').Value = 'Villa Tejúpam De La Unión'
$ws.Range('B605').Value = 'Zimatlán De Álvarez'
$ws.Range('B623').Value = 'Huehuetlán El Grande'
$ws.Range('B627').Value = 'Izúcar De Matamoros'
$ws.Range('B630').Value = 'Palmar De Bravo'
$ws.Range('B640').Value = 'Tepatlaxco De Hidalgo'
$ws.Range('B643').Value = 'Tepexi De Rodríguez'
$ws.Range('B644').Value = 'Tetela De Ocampo'
$ws.Range('B646').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B658').Value = 'Cadereyta De Montes'
$ws.Range('B661').Value = 'Jalpan De Serra'
$ws.Range('B662').Value = 'Landa De Matamoros'
$ws.Range('B664').Value = 'Pinal De Amoles'
$ws.Range('B667').Value = 'San Juan Del Río'
$ws.Range('B678').Value = 'Ciudad Del Maíz'
$ws.Range('B685').Value = 'San Ciro De Acosta'
$ws.Range('B689').Value = 'Santa María Del Río'
$ws.Range('B694').Value = 'Villa De Arriaga'
$ws.Range('B695').Value = 'Villa De Guadalupe'
$ws.Range('B696').Value = 'Villa De Ramos'
$ws.Range('B697').Value = 'Villa De Reyes'
$ws.Range('B740').Value = 'Nacozari De García'
$ws.Range('B773').Value = 'Soto La Marina'
$ws.Range('B779').Value = 'Contla De Juan Cuamatzi'
$ws.Range('B794').Value = 'Boca Del Río'
$ws.Range('B802').Value = 'Cosamaloapan De Carpio'
$ws.Range('B815').Value = 'Lerdo De Tejada'
$ws.Range('B816').Value = 'Martínez De La Torre'
$ws.Range('B829').Value = 'Poza Rica De Hidalgo'
$ws.Range('B853').Value = 'Cañitas De Felipe Pescador'
$ws.Range('B855').Value = 'Concepción Del Oro'
$ws.Range('B868').Value = 'Mezquital Del Oro'
$ws.Range('B873').Value = 'Moyahua De Estrada'
$ws.Range('B874').Value = 'Nochistlán De Mejía'
$ws.Range('B875').Value = 'Noria De Ángeles'
$ws.Range('B885').Value = 'Tlaltenango De Sánchez Román'
$ws.Range('B888').Value = 'Villa De Cos'

# --- Tiny floating point recompute of the Santo Domingo percentage cell ---
$ws.Range('D690').Value = 0.009681227863046043

# --- Drop the trailing footnote / source metadata rows (896-900); ---
# --- the table itself ends at row 894, dimension shrinks to A1:D894 ---
$ws.Rows('896:900').Delete() | Out-Null

